$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header row (row 1), pushing the
# existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Populate the newly-inserted row 2 with slug-style identifiers that
# correspond one-to-one with the header columns in row 1. These values
# allow two columns to be related to build a SKOS hierarchy.
$ws.Range("A2").Value = "nivel-estudios-codigo"
$ws.Range("B2").Value = "personas"
$ws.Range("C2").Value = "municipio-codigo"
$ws.Range("D2").Value = "nivel-estudios"
$ws.Range("E2").Value = "municipio-nombre"
